# Week 5 updates:
# 1) Remove the debugging-screenshot slide (p:sldId id="276" — the
#    "Content Placeholder 3" slide that shows a Python traceback
#    screenshot, currently deck position 19).
# 2) Nudge the picture on the "Nobody  Writes Perfect Code" slide
#    (deck position 2) down slightly.

$p = $ppt.ActivePresentation

# --- 1. Delete the traceback screenshot slide (match by persistent SlideID,
#        since indices shift as slides are added/removed elsewhere) ---
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 276) {
        $s.Delete()
    }
}

# --- 2. Move the picture down on slide 2 ("Nobody  Writes Perfect Code") ---
$s2 = $p.Slides.Item(2)
for ($j = 1; $j -le $s2.Shapes.Count; $j++) {
    $sh = $s2.Shapes.Item($j)
    if ($sh.Name -eq "Picture 2") {
        $sh.Top = 143.75
    }
}
